$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update "Förändrad" date column (C) for rows 2-5 from 2023-09-15 (45184) to 2023-09-16 (45185)
$ws.Range("C2").Value = 45185
$ws.Range("C3").Value = 45185
$ws.Range("C4").Value = 45185
$ws.Range("C5").Value = 45185
